$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 1380990

$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = '1380990'
$ws.Range("N8").ClearFormats()

$ws.Range("D15").Value = 857000

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = 'https://kc-klassavto.ru/auto/changan/cs55_plus/'
$ws.Range("E15").ClearFormats()

$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = '1238900'
$ws.Range("L15").ClearFormats()

$ws.Range("N15").NumberFormat = "@"
$ws.Range("N15").Value = '857000'
$ws.Range("N15").ClearFormats()

$ws.Range("D22").Value = 2039900

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = 'https://kc-klassavto.ru/auto/changan/uni-k/'
$ws.Range("E22").ClearFormats()

$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = '2068900'
$ws.Range("L22").ClearFormats()

$ws.Range("N22").NumberFormat = "@"
$ws.Range("N22").Value = '2039900'
$ws.Range("N22").ClearFormats()

$ws.Range("D24").Value = 1769900

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = 'https://kc-klassavto.ru/auto/changan/uni-t/'
$ws.Range("E24").ClearFormats()

$ws.Range("N24").NumberFormat = "@"
$ws.Range("N24").Value = '1769900'
$ws.Range("N24").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = 'https://kcelitauto.ru/auto/changan/uni-v/'
$ws.Range("E25").ClearFormats()

$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value = '1628900'
$ws.Range("L25").ClearFormats()

$ws.Range("N25").NumberFormat = "@"
$ws.Range("N25").Value = '1726900'
$ws.Range("N25").ClearFormats()

$ws.Range("J26").Value = 2279033

$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = '1389000'
$ws.Range("N28").ClearFormats()

$ws.Range("D30").Value = 1652967

$ws.Range("J30").Value = 1652967

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = 'https://kcelitauto.ru/auto/chery/tiggo_7_pro/'
$ws.Range("E33").ClearFormats()

$ws.Range("L33").NumberFormat = "@"
$ws.Range("L33").Value = '1108900'
$ws.Range("L33").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = 'https://kcelitauto.ru/auto/chery/tiggo_7_pro_max_1/'
$ws.Range("E34").ClearFormats()

$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value = '1304000'
$ws.Range("L34").ClearFormats()

$ws.Range("J35").Value = 1913433

$ws.Range("D41").Value = 2287892

$ws.Range("J41").Value = 2287892

$ws.Range("J42").Value = 3447500

$ws.Range("D46").Value = 759000

$ws.Range("L46").NumberFormat = "@"
$ws.Range("L46").Value = '759000'
$ws.Range("L46").ClearFormats()

$ws.Range("N46").NumberFormat = "@"
$ws.Range("N46").Value = '919900'
$ws.Range("N46").ClearFormats()

$ws.Range("D62").Value = 1539000

$ws.Range("E62").NumberFormat = "@"
$ws.Range("E62").Value = 'https://kc-klassavto.ru/auto/faw/besturn_t77/'
$ws.Range("E62").ClearFormats()

$ws.Range("L62").NumberFormat = "@"
$ws.Range("L62").Value = '1552000'
$ws.Range("L62").ClearFormats()

$ws.Range("N62").NumberFormat = "@"
$ws.Range("N62").Value = '1539000'
$ws.Range("N62").ClearFormats()

$ws.Range("D65").Value = 619300

$ws.Range("E65").NumberFormat = "@"
$ws.Range("E65").Value = 'https://kcelitauto.ru/auto/faw/x40/'
$ws.Range("E65").ClearFormats()

$ws.Range("L65").NumberFormat = "@"
$ws.Range("L65").Value = '619300'
$ws.Range("L65").ClearFormats()

$ws.Range("D68").Value = 1299000

$ws.Range("E68").NumberFormat = "@"
$ws.Range("E68").Value = 'https://kc-klassavto.ru/auto/gac/gs3/'
$ws.Range("E68").ClearFormats()

$ws.Range("N68").NumberFormat = "@"
$ws.Range("N68").Value = '1299000'
$ws.Range("N68").ClearFormats()

$ws.Range("N73").NumberFormat = "@"
$ws.Range("N73").Value = '1999990'
$ws.Range("N73").ClearFormats()

$ws.Range("E74").NumberFormat = "@"
$ws.Range("E74").Value = 'https://kcelitauto.ru/auto/geely/atlas_pro/'
$ws.Range("E74").ClearFormats()

$ws.Range("L74").NumberFormat = "@"
$ws.Range("L74").Value = '1087596'
$ws.Range("L74").ClearFormats()

$ws.Range("D82").Value = 2599000

$ws.Range("E82").NumberFormat = "@"
$ws.Range("E82").Value = 'https://kc-klassavto.ru/auto/geely/monjaro/'
$ws.Range("E82").ClearFormats()

$ws.Range("N82").NumberFormat = "@"
$ws.Range("N82").Value = '2599000'
$ws.Range("N82").ClearFormats()

$ws.Range("D92").Value = 1519000

$ws.Range("E92").NumberFormat = "@"
$ws.Range("E92").Value = 'https://kc-klassavto.ru/auto/haval/h3/'
$ws.Range("E92").ClearFormats()

$ws.Range("N92").NumberFormat = "@"
$ws.Range("N92").Value = '1519000'
$ws.Range("N92").ClearFormats()

$ws.Range("D98").Value = 689900

$ws.Range("E98").NumberFormat = "@"
$ws.Range("E98").Value = 'https://kc-klassavto.ru/auto/haval/jolion_1/'
$ws.Range("E98").ClearFormats()

$ws.Range("N98").NumberFormat = "@"
$ws.Range("N98").Value = '689900'
$ws.Range("N98").ClearFormats()

$ws.Range("E101").NumberFormat = "@"
$ws.Range("E101").Value = 'https://kcelitauto.ru/auto/hyundai/creta/'
$ws.Range("E101").ClearFormats()

$ws.Range("L101").NumberFormat = "@"
$ws.Range("L101").Value = '676000'
$ws.Range("L101").ClearFormats()

$ws.Range("E106").NumberFormat = "@"
$ws.Range("E106").Value = 'https://kcelitauto.ru/auto/hyundai/solaris/'
$ws.Range("E106").ClearFormats()

$ws.Range("L106").NumberFormat = "@"
$ws.Range("L106").Value = '465000'
$ws.Range("L106").ClearFormats()

$ws.Range("N107").NumberFormat = "@"
$ws.Range("N107").Value = '1492000'
$ws.Range("N107").ClearFormats()

$ws.Range("E110").NumberFormat = "@"
$ws.Range("E110").Value = 'https://kcelitauto.ru/auto/jac/j7/'
$ws.Range("E110").ClearFormats()

$ws.Range("L110").NumberFormat = "@"
$ws.Range("L110").Value = '718000'
$ws.Range("L110").ClearFormats()

$ws.Range("N110").NumberFormat = "@"
$ws.Range("N110").Value = '899000'
$ws.Range("N110").ClearFormats()

$ws.Range("J120").Value = 2329900

$ws.Range("L146").NumberFormat = "@"
$ws.Range("L146").Value = '425900'
$ws.Range("L146").ClearFormats()

$ws.Range("N146").NumberFormat = "@"
$ws.Range("N146").Value = '462010'
$ws.Range("N146").ClearFormats()

$ws.Range("L147").NumberFormat = "@"
$ws.Range("L147").Value = '272300'
$ws.Range("L147").ClearFormats()

$ws.Range("E148").NumberFormat = "@"
$ws.Range("E148").Value = 'https://kcelitauto.ru/auto/lada/granta_liftback/'
$ws.Range("E148").ClearFormats()

$ws.Range("L148").NumberFormat = "@"
$ws.Range("L148").Value = '268500'
$ws.Range("L148").ClearFormats()

$ws.Range("N148").NumberFormat = "@"
$ws.Range("N148").Value = '375610'
$ws.Range("N148").ClearFormats()

$ws.Range("E150").NumberFormat = "@"
$ws.Range("E150").Value = 'https://kcelitauto.ru/auto/lada/granta_sedan/'
$ws.Range("E150").ClearFormats()

$ws.Range("L150").NumberFormat = "@"
$ws.Range("L150").Value = '257900'
$ws.Range("L150").ClearFormats()

$ws.Range("N150").NumberFormat = "@"
$ws.Range("N150").Value = '358510'
$ws.Range("N150").ClearFormats()

$ws.Range("L151").NumberFormat = "@"
$ws.Range("L151").Value = '425900'
$ws.Range("L151").ClearFormats()

$ws.Range("L157").NumberFormat = "@"
$ws.Range("L157").Value = '278900'
$ws.Range("L157").ClearFormats()

$ws.Range("E158").NumberFormat = "@"
$ws.Range("E158").Value = 'https://kcelitauto.ru/auto/lada/largus/'
$ws.Range("E158").ClearFormats()

$ws.Range("L158").NumberFormat = "@"
$ws.Range("L158").Value = '348900'
$ws.Range("L158").ClearFormats()

$ws.Range("E159").NumberFormat = "@"
$ws.Range("E159").Value = 'https://kcelitauto.ru/auto/lada/largus_cross/'
$ws.Range("E159").ClearFormats()

$ws.Range("L159").NumberFormat = "@"
$ws.Range("L159").Value = '384300'
$ws.Range("L159").ClearFormats()

$ws.Range("D160").Value = 699310

$ws.Range("N160").NumberFormat = "@"
$ws.Range("N160").Value = '699310'
$ws.Range("N160").ClearFormats()

$ws.Range("D161").Value = 540310

$ws.Range("N161").NumberFormat = "@"
$ws.Range("N161").Value = '540310'
$ws.Range("N161").ClearFormats()

$ws.Range("L167").NumberFormat = "@"
$ws.Range("L167").Value = '380200'
$ws.Range("L167").ClearFormats()

$ws.Range("D169").Value = 444500

$ws.Range("E169").NumberFormat = "@"
$ws.Range("E169").Value = 'https://kcelitauto.ru/auto/lada/vesta_cross/'
$ws.Range("E169").ClearFormats()

$ws.Range("L169").NumberFormat = "@"
$ws.Range("L169").Value = '444500'
$ws.Range("L169").ClearFormats()

$ws.Range("D172").Value = 433000

$ws.Range("E172").NumberFormat = "@"
$ws.Range("E172").Value = 'https://kcelitauto.ru/auto/lada/vesta_sw_cross/'
$ws.Range("E172").ClearFormats()

$ws.Range("L172").NumberFormat = "@"
$ws.Range("L172").Value = '433000'
$ws.Range("L172").ClearFormats()

$ws.Range("D173").Value = 390000

$ws.Range("E173").NumberFormat = "@"
$ws.Range("E173").Value = 'https://kcelitauto.ru/auto/lada/vesta_sw/'
$ws.Range("E173").ClearFormats()

$ws.Range("L173").NumberFormat = "@"
$ws.Range("L173").Value = '390000'
$ws.Range("L173").ClearFormats()

$ws.Range("D177").Value = 380900

$ws.Range("L177").NumberFormat = "@"
$ws.Range("L177").Value = '380900'
$ws.Range("L177").ClearFormats()

$ws.Range("L180").NumberFormat = "@"
$ws.Range("L180").Value = '385900'
$ws.Range("L180").ClearFormats()

$ws.Range("D192").Value = 1249000

$ws.Range("E192").NumberFormat = "@"
$ws.Range("E192").Value = 'https://kc-klassavto.ru/auto/moskvich/6/'
$ws.Range("E192").ClearFormats()

$ws.Range("N192").NumberFormat = "@"
$ws.Range("N192").Value = '1249000'
$ws.Range("N192").ClearFormats()

$ws.Range("L193").NumberFormat = "@"
$ws.Range("L193").Value = '897000'
$ws.Range("L193").ClearFormats()

$ws.Range("E195").NumberFormat = "@"
$ws.Range("E195").Value = 'https://kcelitauto.ru/auto/nissan/x-trail/'
$ws.Range("E195").ClearFormats()

$ws.Range("L195").NumberFormat = "@"
$ws.Range("L195").Value = '1134000'
$ws.Range("L195").ClearFormats()

$ws.Range("N196").NumberFormat = "@"
$ws.Range("N196").Value = '1318100'
$ws.Range("N196").ClearFormats()

$ws.Range("N203").NumberFormat = "@"
$ws.Range("N203").Value = '938000'
$ws.Range("N203").ClearFormats()

$ws.Range("E204").NumberFormat = "@"
$ws.Range("E204").Value = 'https://kcelitauto.ru/auto/renault/duster/'
$ws.Range("E204").ClearFormats()

$ws.Range("L204").NumberFormat = "@"
$ws.Range("L204").Value = '475000'
$ws.Range("L204").ClearFormats()

$ws.Range("L205").NumberFormat = "@"
$ws.Range("L205").Value = '713000'
$ws.Range("L205").ClearFormats()

$ws.Range("L206").NumberFormat = "@"
$ws.Range("L206").Value = '395000'
$ws.Range("L206").ClearFormats()

$ws.Range("L207").NumberFormat = "@"
$ws.Range("L207").Value = '503990'
$ws.Range("L207").ClearFormats()

$ws.Range("L209").NumberFormat = "@"
$ws.Range("L209").Value = '395000'
$ws.Range("L209").ClearFormats()

$ws.Range("E210").NumberFormat = "@"
$ws.Range("E210").Value = 'https://kcelitauto.ru/auto/renault/sandero_stepway/'
$ws.Range("E210").ClearFormats()

$ws.Range("L210").NumberFormat = "@"
$ws.Range("L210").Value = '399000'
$ws.Range("L210").ClearFormats()

$ws.Range("C211").NumberFormat = "@"
$ws.Range("C211").Value = 'Sandero Stepway City'
$ws.Range("C211").ClearFormats()

$ws.Range("L214").NumberFormat = "@"
$ws.Range("L214").Value = '844800'
$ws.Range("L214").ClearFormats()

$ws.Range("L216").NumberFormat = "@"
$ws.Range("L216").Value = '425000'
$ws.Range("L216").ClearFormats()

$ws.Range("D219").Value = 1499000

$ws.Range("E219").NumberFormat = "@"
$ws.Range("E219").Value = 'https://kc-klassavto.ru/auto/solaris/hc/'
$ws.Range("E219").ClearFormats()

$ws.Range("N219").NumberFormat = "@"
$ws.Range("N219").Value = '1499000'
$ws.Range("N219").ClearFormats()

$ws.Range("D220").Value = 1199000

$ws.Range("E220").NumberFormat = "@"
$ws.Range("E220").Value = 'https://kc-klassavto.ru/auto/solaris/hs/'
$ws.Range("E220").ClearFormats()

$ws.Range("N220").NumberFormat = "@"
$ws.Range("N220").Value = '1199000'
$ws.Range("N220").ClearFormats()

$ws.Range("D221").Value = 1209000

$ws.Range("E221").NumberFormat = "@"
$ws.Range("E221").Value = 'https://kc-klassavto.ru/auto/solaris/krs/'
$ws.Range("E221").ClearFormats()

$ws.Range("N221").NumberFormat = "@"
$ws.Range("N221").Value = '1209000'
$ws.Range("N221").ClearFormats()

$ws.Range("D222").Value = 1239000

$ws.Range("E222").NumberFormat = "@"
$ws.Range("E222").Value = 'https://kc-klassavto.ru/auto/solaris/krx/'
$ws.Range("E222").ClearFormats()

$ws.Range("N222").NumberFormat = "@"
$ws.Range("N222").Value = '1239000'
$ws.Range("N222").ClearFormats()

$ws.Range("E234").NumberFormat = "@"
$ws.Range("E234").Value = 'https://kcelitauto.ru/auto/uaz/patriot/'
$ws.Range("E234").ClearFormats()

$ws.Range("L234").NumberFormat = "@"
$ws.Range("L234").Value = '429380'
$ws.Range("L234").ClearFormats()

$ws.Range("D236").Value = 621000

$ws.Range("E236").NumberFormat = "@"
$ws.Range("E236").Value = 'https://kcelitauto.ru/auto/uaz/patriot_pickup/'
$ws.Range("E236").ClearFormats()

$ws.Range("L236").NumberFormat = "@"
$ws.Range("L236").Value = '621000'
$ws.Range("L236").ClearFormats()

$ws.Range("E241").NumberFormat = "@"
$ws.Range("E241").Value = 'https://kcelitauto.ru/auto/volkswagen/polo/'
$ws.Range("E241").ClearFormats()

$ws.Range("L241").NumberFormat = "@"
$ws.Range("L241").Value = '519300'
$ws.Range("L241").ClearFormats()

$ws.Range("E243").NumberFormat = "@"
$ws.Range("E243").Value = 'https://kcelitauto.ru/auto/volkswagen/tiguan/'
$ws.Range("E243").ClearFormats()

$ws.Range("L243").NumberFormat = "@"
$ws.Range("L243").Value = '1288900'
$ws.Range("L243").ClearFormats()

$ws.Range("J244").Value = 2130560
